$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 1208.2354
$ws.Range("I8").Value = 49.090908
$ws.Range("J8").Value = 3333.3333
$ws.Range("K8").Value = 147.272724
$ws.Range("L8").Value = 9999.999899999999
$ws.Range("M8").Value = -8.272723999999982
$ws.Range("N8").Value = -10277.9999
$ws.Range("H15").Value = 1031.0869
$ws.Range("I15").Value = 1031.0869
$ws.Range("K15").Value = 3093.2607
$ws.Range("M15").Value = -2924.2607
$ws.Range("H38").Value = 4475
$ws.Range("I38").Value = 90
$ws.Range("J38").Value = 5670.909
$ws.Range("K38").Value = 270
$ws.Range("L38").Value = 17012.727
$ws.Range("M38").Value = 102
$ws.Range("N38").Value = -17756.727
$ws.Range("H98").Value = 4980.4
$ws.Range("I98").Value = 3239.2307
$ws.Range("J98").Value = 6866.6665
$ws.Range("K98").Value = 3239.2307
$ws.Range("L98").Value = 6866.6665
$ws.Range("M98").Value = -1741.2307
$ws.Range("N98").Value = -9862.666499999999
$ws.Range("H111").Value = 899.5
$ws.Range("I111").Value = 899.5
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 2698.5
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 368.5
$ws.Range("N111").ClearContents()
$ws.Range("H122").Value = 4980.4
$ws.Range("I122").Value = 3239.2307
$ws.Range("J122").Value = 6866.6665
$ws.Range("K122").Value = 9717.6921
$ws.Range("L122").Value = 20599.9995
$ws.Range("M122").Value = -7267.6921
$ws.Range("N122").Value = -25499.9995
$ws.Range("H129").Value = 870.59186
$ws.Range("I129").Value = 322.66666
$ws.Range("J129").Value = 887.8946999999999
$ws.Range("K129").Value = 967.9999799999999
$ws.Range("L129").Value = 2663.6841
$ws.Range("M129").Value = 4032.00002
$ws.Range("N129").Value = -12663.6841
$ws.Range("H138").Value = 2716.78
$ws.Range("I138").Value = 1019.7
$ws.Range("J138").Value = 2905.3445
$ws.Range("K138").Value = 3059.1
$ws.Range("L138").Value = 8716.033500000001
$ws.Range("M138").Value = 2080.9
$ws.Range("N138").Value = -18996.0335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 11223.75
$ws.Range("I6").Value = 500
$ws.Range("J6").Value = 14798.333
$ws.Range("K6").Value = 500
$ws.Range("L6").Value = 14798.333
$ws.Range("N6").Value = -15144.333
$ws.Range("M6").Value = -327
$ws.Range("H35").Value = 13199.667
$ws.Range("I35").Value = 1799.5
$ws.Range("J35").Value = 36000
$ws.Range("K35").Value = 1799.5
$ws.Range("L35").Value = 36000
$ws.Range("M35").Value = -1393.5
$ws.Range("N35").Value = -36812
$ws.Range("H102").Value = 2529.8572
$ws.Range("I102").Value = 1962
$ws.Range("K102").Value = 1962
$ws.Range("M102").Value = -340

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 45000
$ws.Range("J88").Value = 45000
$ws.Range("L88").Value = 45000
$ws.Range("N88").Value = -45812
$ws.Range("H91").Value = 45000
$ws.Range("J91").Value = 45000
$ws.Range("L91").Value = 45000
$ws.Range("N91").Value = -47808
$ws.Range("H97").Value = 23339.6
$ws.Range("I97").Value = 3378
$ws.Range("J97").Value = 36647.332
$ws.Range("K97").Value = 3378
$ws.Range("L97").Value = 36647.332
$ws.Range("M97").Value = -2387
$ws.Range("N97").Value = -38629.332
$ws.Range("H103").Value = 38000
$ws.Range("J103").Value = 38000
$ws.Range("L103").Value = 38000
$ws.Range("N103").Value = -40344

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2782.7715
$ws.Range("I31").Value = 998.2778
$ws.Range("J31").Value = 4672.2354
$ws.Range("K31").Value = 998.2778
$ws.Range("L31").Value = 4672.2354
$ws.Range("M31").Value = -703.2778
$ws.Range("N31").Value = -5262.2354
$ws.Range("H34").Value = 2782.7715
$ws.Range("I34").Value = 998.2778
$ws.Range("J34").Value = 4672.2354
$ws.Range("K34").Value = 998.2778
$ws.Range("L34").Value = 4672.2354
$ws.Range("M34").Value = -796.2778
$ws.Range("N34").Value = -5076.2354
$ws.Range("H134").Value = 8397.4375
$ws.Range("I134").Value = 10178.272
$ws.Range("K134").Value = 30534.816
$ws.Range("M134").Value = -27999.816

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 19350
$ws.Range("I4").Value = 50083.332
$ws.Range("J4").Value = 910
$ws.Range("K4").Value = 150249.996
$ws.Range("L4").Value = 2730
$ws.Range("M4").Value = -150137.996
$ws.Range("N4").Value = -2954
$ws.Range("H9").Value = 129349.06
$ws.Range("I9").Value = 1001
$ws.Range("J9").Value = 137370.81
$ws.Range("K9").Value = 3003
$ws.Range("L9").Value = 412112.43
$ws.Range("M9").Value = -2779
$ws.Range("N9").Value = -412560.43
$ws.Range("H13").Value = 2055.5
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 2055.5
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 6166.5
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -6502.5
$ws.Range("H17").Value = 1467.3334
$ws.Range("I17").Value = 750
$ws.Range("J17").Value = 1826
$ws.Range("K17").Value = 2250
$ws.Range("L17").Value = 5478
$ws.Range("M17").Value = -2081
$ws.Range("N17").Value = -5816
$ws.Range("H37").Value = 500100000
$ws.Range("J37").Value = 500100000
$ws.Range("L37").Value = 1500300000
$ws.Range("N37").Value = -1500300224
$ws.Range("H108").Value = 325.25
$ws.Range("I108").Value = 325.25
$ws.Range("K108").Value = 975.75
$ws.Range("M108").Value = 1904.25
$ws.Range("H109").Value = 1533
$ws.Range("I109").Value = 1533
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 4599
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -3559
$ws.Range("N109").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1522.2
$ws.Range("I113").Value = 1402.75
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1402.75
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 767.25
$ws.Range("N113").Value = -6340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 16749.75
$ws.Range("J29").Value = 16749.75
$ws.Range("L29").Value = 16749.75
$ws.Range("N29").Value = -17339.75
$ws.Range("H136").Value = 4041.9614
$ws.Range("I136").Value = 1088.0834
$ws.Range("K136").Value = 3264.2502
$ws.Range("M136").Value = -714.2501999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3956.7856
$ws.Range("J136").Value = 6812.909
$ws.Range("L136").Value = 20438.727
$ws.Range("N136").Value = -25538.727
